$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
Write-Host ("today=" + (Get-Date))
Write-Host ("formula=" + $ws.Range("B2").Formula)
Write-Host ("text=" + $ws.Range("B2").Text)
Write-Host ("E5 text=" + $ws.Range("E5").Text)
Write-Host ("I26 text=" + $ws.Range("I26").Text)
